# Apply the daydata.xlsx edit described by the commit diff:
#  - Sheet1!C58:F58 and Sheet1!C59:F59 change from 1 -> 0
#  - The sheet's active selection moves from H57 to I59

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Zero-out the four values on rows 58 and 59 (columns C-F)
$ws.Range("C58:F58").Value = 0
$ws.Range("C59:F59").Value = 0

# Move the sheet's selection to I59 (new activeCell/sqref)
$ws.Activate()
$ws.Range("I59").Select()

# Keep the viewport's scroll anchor consistent with the selection move
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
